# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's rush/pass yardage-per-play logs (OFF/DEF)
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$r = $ydsWs.Range("B2")
$r.Value = $r.Value() + " 4 8 3 1 9 2 0 1 7 5 3 1 0 7 3 4 0 5 4 12 4"

$r = $ydsWs.Range("B3")
$r.Value = $r.Value() + " 4 10 5 17 23 3 26 16 11 6 -2 6 8 4 11 17 19 8 12 15 6 12 56 2"

$r = $ydsWs.Range("C2")
$r.Value = $r.Value() + " 7 0 0 8 9 -1 0 0 3 1 0 5 10 3 2 4 3 -1 4 5 23 5 1 -1"

$r = $ydsWs.Range("C3")
$r.Value = $r.Value() + " -8 14 3 7 -3 11 4 0 12 0 24 18 -1 13 2 42 12 18 18 13 6 4"

# ---------------------------------------------------------------------------
# OFF sheet: Week 16 box-score totals
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 209
$offWs.Range("D2").Value = 13
$offWs.Range("F2").Value = 62
$offWs.Range("G2").Value = 57
$offWs.Range("J2").Value = 30
$offWs.Range("N2").Value = 13
$offWs.Range("O2").Value = 18
$offWs.Range("P2").Value = 10

$offWs.Range("B3").Value = 13
$offWs.Range("C3").Value = 146
$offWs.Range("E3").Value = 31
$offWs.Range("F3").Value = 89
$offWs.Range("G3").Value = 24
$offWs.Range("H3").Value = 29
$offWs.Range("I3").Value = 50
$offWs.Range("J3").Value = 35
$offWs.Range("L3").Value = 237
$offWs.Range("M3").Value = 159
$offWs.Range("Q3").Value = 482

# ---------------------------------------------------------------------------
# DEF sheet: Week 16 box-score totals
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("B2").Value = 5
$defWs.Range("C2").Value = 177
$defWs.Range("E2").Value = 15
$defWs.Range("F2").Value = 57
$defWs.Range("G2").Value = 47
$defWs.Range("H2").Value = 5
$defWs.Range("J2").Value = 24
$defWs.Range("N2").Value = 27

$defWs.Range("B3").Value = 10
$defWs.Range("C3").Value = 162
$defWs.Range("D3").Value = 5
$defWs.Range("E3").Value = 34
$defWs.Range("F3").Value = 89
$defWs.Range("H3").Value = 29
$defWs.Range("I3").Value = 56
$defWs.Range("J3").Value = 55
$defWs.Range("L3").Value = 265
$defWs.Range("M3").Value = 184
$defWs.Range("Q3").Value = 502

# ---------------------------------------------------------------------------
# ST sheet: Week 16 KO/PT totals plus appended KO/PT distance & return logs
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value = 77
$stWs.Range("D2").Value = 47
$stWs.Range("F2").Value = 582
$stWs.Range("G2").Value = 567
$stWs.Range("L2").Value = 161
$stWs.Range("M2").Value = 126
$stWs.Range("B3").Value = 34

$r = $stWs.Range("B4")
$r.Value = $r.Value() + " 63 60"

$r = $stWs.Range("B5")
$r.Value = $r.Value() + " 37 20"

$r = $stWs.Range("D3")
$r.Value = $r.Value() + " 30 40"

$r = $stWs.Range("D4")
$r.Value = $r.Value() + " 0 0"

$r = $stWs.Range("D5")
$r.Value = $r.Value() + " 0 0 12 0 0"

# ---------------------------------------------------------------------------
# TURNS sheet: Week 16 turnover totals
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B3").Value = 5
$turnsWs.Range("D3").Value = 10
$turnsWs.Range("E3").Value = 11

# ---------------------------------------------------------------------------
# PEN sheet: Week 16 penalty totals
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 17
$penWs.Range("B3").Value = 16
$penWs.Range("D4").Value = 17
